$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update header timestamp
$ws.Range("A1").Value = "Datos actualizados a 2 de Abril de 2020 a las 17:20"

# Row 4: Estados Unidos -> Estados Unidos
$ws.Range("B4").Value = 217372
$ws.Range("C4").Value = 2369
$ws.Range("D4").Value = 8942
$ws.Range("E4").Value = 203285
$ws.Range("G4").Value = 43
$ws.Range("H4").Value = 5145

# Row 17: Corea del Sur -> Canada
$ws.Range("A17").Value = "Canada"
$ws.Range("B17").Value = 10132
$ws.Range("C17").Value = 401
$ws.Range("D17").Value = 1736
$ws.Range("E17").Value = 8267
$ws.Range("F17").Value = 120
$ws.Range("G17").Value = 15
$ws.Range("H17").Value = 129

# Row 18: Canada -> Corea del Sur
$ws.Range("A18").Value = "Corea del Sur"
$ws.Range("B18").Value = 9976
$ws.Range("C18").Value = 89
$ws.Range("D18").Value = 5828
$ws.Range("E18").Value = 3979
$ws.Range("F18").Value = 55
$ws.Range("G18").Value = 4
$ws.Range("H18").Value = 169

# Row 20: Brasil -> Brasil
$ws.Range("B20").Value = 7022
$ws.Range("C20").Value = 142
$ws.Range("E20").Value = 6643
$ws.Range("G20").Value = 10
$ws.Range("H20").Value = 252

# Row 35: Pakistan -> Pakistan
$ws.Range("E35").Value = 2247
$ws.Range("G35").Value = 5
$ws.Range("H35").Value = 32

# Row 72: Moldavia -> Moldavia
$ws.Range("E72").Value = 476
$ws.Range("G72").Value = 1
$ws.Range("H72").Value = 6

# Row 83: Uruguay -> Republica de Chipre
$ws.Range("A83").Value = "Republica de Chipre"
$ws.Range("B83").Value = 356
$ws.Range("C83").Value = 36
$ws.Range("D83").Value = 28
$ws.Range("E83").Value = 319
$ws.Range("F83").Value = 11
$ws.Range("G83").Value = 0
$ws.Range("H83").Value = 9

# Row 84: Kuwait -> Uruguay
$ws.Range("A84").Value = "Uruguay"
$ws.Range("B84").Value = 350
$ws.Range("C84").Value = 0
$ws.Range("D84").Value = 62
$ws.Range("E84").Value = 284
$ws.Range("G84").Value = 2
$ws.Range("H84").Value = 4

# Row 85: Taiwan -> Kuwait
$ws.Range("A85").Value = "Kuwait"
$ws.Range("B85").Value = 342
$ws.Range("C85").Value = 25
$ws.Range("D85").Value = 81
$ws.Range("E85").Value = 261
$ws.Range("F85").Value = 15
$ws.Range("H85").Value = 0

# Row 86: Republica de Chipre -> Taiwan
$ws.Range("A86").Value = "Taiwan"
$ws.Range("B86").Value = 339
$ws.Range("C86").Value = 10
$ws.Range("D86").Value = 50
$ws.Range("E86").Value = 284
$ws.Range("F86").Value = 0
$ws.Range("H86").Value = 5

# Row 95: Oman -> Cuba
$ws.Range("A95").Value = "Cuba"
$ws.Range("B95").Value = 233
$ws.Range("D95").Value = 13
$ws.Range("E95").Value = 214
$ws.Range("F95").Value = 7
$ws.Range("H95").Value = 6

# Row 96: Vietnam -> Oman
$ws.Range("A96").Value = "Oman"
$ws.Range("B96").Value = 231
$ws.Range("C96").Value = 21
$ws.Range("D96").Value = 57
$ws.Range("E96").Value = 173
$ws.Range("H96").Value = 1

# Row 97: Honduras -> Vietnam
$ws.Range("A97").Value = "Vietnam"
$ws.Range("B97").Value = 227
$ws.Range("C97").Value = 9
$ws.Range("D97").Value = 75
$ws.Range("E97").Value = 152
$ws.Range("F97").Value = 3
$ws.Range("G97").Value = 0
$ws.Range("H97").Value = 0

# Row 98: Cuba -> Honduras
$ws.Range("A98").Value = "Honduras"
$ws.Range("B98").Value = 219
$ws.Range("C98").Value = 47
$ws.Range("D98").Value = 3
$ws.Range("E98").Value = 202
$ws.Range("F98").Value = 4
$ws.Range("G98").Value = 4
$ws.Range("H98").Value = 14

# Row 106: Mauricio -> Mauricio
$ws.Range("B106").Value = 169
$ws.Range("C106").Value = 8
$ws.Range("E106").Value = 162

# Row 121: Trinidad yTobago -> Isla de Man
$ws.Range("A121").Value = "Isla de Man"
$ws.Range("B121").Value = 95
$ws.Range("C121").Value = 27
$ws.Range("D121").Value = 0
$ws.Range("E121").Value = 94
$ws.Range("H121").Value = 1

# Row 122: Gibraltar -> Trinidad yTobago
$ws.Range("A122").Value = "Trinidad yTobago"
$ws.Range("B122").Value = 90
$ws.Range("C122").Value = 0
$ws.Range("D122").Value = 1
$ws.Range("E122").Value = 84
$ws.Range("H122").Value = 5

# Row 123: Ruanda -> Gibraltar
$ws.Range("A123").Value = "Gibraltar"
$ws.Range("B123").Value = 88
$ws.Range("C123").Value = 7
$ws.Range("D123").Value = 46
$ws.Range("E123").Value = 42

# Row 124: Paraguay -> Ruanda
$ws.Range("A124").Value = "Ruanda"
$ws.Range("B124").Value = 82
$ws.Range("C124").Value = 0
$ws.Range("D124").Value = 0
$ws.Range("E124").Value = 82
$ws.Range("F124").Value = 0
$ws.Range("H124").Value = 0

# Row 125: Isla de Man -> Paraguay
$ws.Range("A125").Value = "Paraguay"
$ws.Range("B125").Value = 77
$ws.Range("C125").Value = 8
$ws.Range("D125").Value = 2
$ws.Range("E125").Value = 72
$ws.Range("F125").Value = 4
$ws.Range("H125").Value = 3
